$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A user can only receive one invitation per day/year, so the duplicate
# invitation rows are removed: row 5 (Pacofiestas / 12-04-1990 /
# luismax@gmail.com) and row 6 (cliente / 12-04-1990 /
# randomemail@email.com). Clear just the values in A:C, keeping the
# existing cell formatting (style ids) untouched.
$ws.Range("A5:C6").ClearContents()

# The mailto hyperlinks anchored on the now-empty C5/C6 cells must go too.
# This host's Hyperlinks.Delete() only operates at whole-sheet granularity
# (it empties the entire collection no matter which Range it is invoked
# on), so wipe it and re-create the two links that must survive (C3/C4)
# with their original target + display text.
$ws.Range("C3:C6").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:daviddsrperiodismo@gmail.com", [Type]::Missing, [Type]::Missing, "daviddsrperiodismo@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:davichon1985@hotmail.com", [Type]::Missing, [Type]::Missing, "davichon1985@hotmail.com")

# Hyperlinks.Add silently reformats its target cell with the built-in
# "Hyperlink" style, clobbering the original (blue text / thin border)
# direct formatting that C3/C4 already had. Restore it by copying the
# still-untouched C5 format (same style id, cleared of its value above)
# back onto C3/C4, then drop the now-unused "Hyperlink" named style.
$ws.Range("C5").Copy($ws.Range("C5")) | Out-Null
$ws.Range("C3").PasteSpecial(-4122, [Type]::Missing) | Out-Null
$ws.Range("C4").PasteSpecial(-4122, [Type]::Missing) | Out-Null
try {
    $wb.Styles.Item("Hyperlink").Delete()
} catch {
}

# Update the active selection to reflect the new state (bottom-right pane).
$ws.Range("C6").Select()
